$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Range("B12").Value = "OPQA-2015||OPQA-3650"
